# Update header text: Volume/Number and report date range
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared string rich-text runs) ---
# A8: "Volume 32   Number  51" -> "...52"
$ws.Range("A8").Characters(21,2).Text = "52"
# C9: "Report Covering the Week  12/15/2025  Through  12/21/2025" -> new dates
$ws.Range("C9").Characters(27,10).Text = "12/22/2025"
$ws.Range("C9").Characters(48,10).Text = "12/28/2025"

# --- Numeric cell updates ---
$ws.Range("C16").Value2 = 1
$ws.Range("G16").Value2 = 2
$ws.Range("H16").Value2 = 300
$ws.Range("I16").Value2 = 84
$ws.Range("K16").Value2 = -5.617977528089
$ws.Range("L16").Value2 = 3.703703703703
$ws.Range("C17").Value2 = 3
$ws.Range("D17").Value2 = 6
$ws.Range("E17").Value2 = -50
$ws.Range("F17").Value2 = 18
$ws.Range("G17").Value2 = 23
$ws.Range("H17").Value2 = -21.739130434782
$ws.Range("I17").Value2 = 247
$ws.Range("J17").Value2 = 211
$ws.Range("K17").Value2 = 17.061611374407
$ws.Range("L17").Value2 = 1.646090534979
$ws.Range("C18").Value2 = 2
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = -50
$ws.Range("G18").Value2 = 8
$ws.Range("H18").Value2 = 0
$ws.Range("I18").Value2 = 79
$ws.Range("J18").Value2 = 103
$ws.Range("K18").Value2 = -23.300970873786
$ws.Range("L18").Value2 = 9.722222222222
$ws.Range("C19").Value2 = 7
$ws.Range("D19").Value2 = 6
$ws.Range("E19").Value2 = 16.666666666666
$ws.Range("F19").Value2 = 26
$ws.Range("G19").Value2 = 27
$ws.Range("H19").Value2 = -3.703703703703
$ws.Range("I19").Value2 = 418
$ws.Range("J19").Value2 = 357
$ws.Range("K19").Value2 = 17.086834733893
$ws.Range("L19").Value2 = -0.712589073634
$ws.Range("D20").Value2 = 1
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 6
$ws.Range("G20").Value2 = 5
$ws.Range("H20").Value2 = 20
$ws.Range("I20").Value2 = 76
$ws.Range("J20").Value2 = 90
$ws.Range("K20").Value2 = -15.555555555555
$ws.Range("L20").Value2 = -31.531531531531
$ws.Range("C21").Value2 = 14
$ws.Range("D21").Value2 = 17
$ws.Range("E21").Value2 = -17.647058823529
$ws.Range("F21").Value2 = 67
$ws.Range("G21").Value2 = 66
$ws.Range("H21").Value2 = 1.515151515151
$ws.Range("I21").Value2 = 929
$ws.Range("J21").Value2 = 872
$ws.Range("K21").Value2 = 6.536697247706
$ws.Range("L21").Value2 = -1.588983050847
$ws.Range("D23").Value2 = 1
$ws.Range("E23").Value2 = -100
$ws.Range("G23").Value2 = 1
$ws.Range("H23").Value2 = -100
$ws.Range("J23").Value2 = 18
$ws.Range("K23").Value2 = 5.555555555555
$ws.Range("C24").Value2 = 20
$ws.Range("D24").Value2 = 32
$ws.Range("E24").Value2 = -37.5
$ws.Range("F24").Value2 = 82
$ws.Range("H24").Value2 = -47.435897435897
$ws.Range("I24").Value2 = 1312
$ws.Range("J24").Value2 = 1487
$ws.Range("K24").Value2 = -11.768661735037
$ws.Range("L24").Value2 = -10.990502035278
$ws.Range("C25").Value2 = 16
$ws.Range("D25").Value2 = 28
$ws.Range("E25").Value2 = -42.857142857142
$ws.Range("F25").Value2 = 58
$ws.Range("G25").Value2 = 112
$ws.Range("H25").Value2 = -48.214285714285
$ws.Range("I25").Value2 = 909
$ws.Range("J25").Value2 = 1018
$ws.Range("K25").Value2 = -10.707269155206
$ws.Range("L25").Value2 = 7.446808510638
$ws.Range("D26").Value2 = 2
$ws.Range("E26").Value2 = 400
$ws.Range("F26").Value2 = 47
$ws.Range("G26").Value2 = 26
$ws.Range("H26").Value2 = 80.76923076923
$ws.Range("I26").Value2 = 575
$ws.Range("J26").Value2 = 563
$ws.Range("K26").Value2 = 2.131438721136
$ws.Range("L26").Value2 = 12.085769980506
$ws.Range("C28").Value2 = 1
$ws.Range("I28").Value2 = 59
$ws.Range("K28").Value2 = -1.666666666666
$ws.Range("L28").Value2 = -4.838709677419
$ws.Range("C33").Value2 = 1
$ws.Range("F33").Value2 = 1
$ws.Range("I33").Value2 = 5
$ws.Range("K33").Value2 = 66.666666666666
$ws.Range("L33").Value2 = 400

# --- Text ("N/A"-style) cell updates: values come from shared strings "0" (idx 20) and "***.*" (idx 21) ---
# Use an apostrophe prefix to force text interpretation, then copy the number format
# from a cell that already carries the correct "N/A" style (style index 13) so no new
# style entry is created and the cell keeps matching the original formatting exactly.
$ws.Range("G14").Value2 = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("G14").PasteSpecial(-4122) | Out-Null
$ws.Range("H14").Value2 = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("H14").PasteSpecial(-4122) | Out-Null
$ws.Range("D15").Value2 = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").Value2 = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").Value2 = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value2 = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("D28").Value2 = "'0"
$ws.Range("C14").Copy() | Out-Null
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value2 = "'***.*"
$ws.Range("E14").Copy() | Out-Null
$ws.Range("E28").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
